$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-7 (weekly price refresh) ---

# Row 2
$ws.Cells.Item(2,4).Value = 44761   # D2 Fecha
$ws.Cells.Item(2,10).Value = 200    # J2 Volumen
$ws.Cells.Item(2,11).Value = 700    # K2 Precio minimo
$ws.Cells.Item(2,12).Value = 800    # L2 Precio maximo
$ws.Cells.Item(2,13).Value = 750    # M2 Precio promedio ponderado
$ws.Cells.Item(2,16).Value = 750    # P2 Precio $/Kg

# Row 3
$ws.Cells.Item(3,4).Value = 44761   # D3 Fecha
$ws.Cells.Item(3,9).Value = "Segunda" # I3 Calidad
$ws.Cells.Item(3,10).Value = 150    # J3 Volumen
$ws.Cells.Item(3,11).Value = 600    # K3 Precio minimo
$ws.Cells.Item(3,12).Value = 600    # L3 Precio maximo
$ws.Cells.Item(3,13).Value = 600    # M3 Precio promedio ponderado
$ws.Cells.Item(3,16).Value = 600    # P3 Precio $/Kg

# Row 4
$ws.Cells.Item(4,4).Value = 44610   # D4 Fecha
$ws.Cells.Item(4,10).Value = 100    # J4 Volumen
$ws.Cells.Item(4,11).Value = 600    # K4 Precio minimo
$ws.Cells.Item(4,12).Value = 650    # L4 Precio maximo
$ws.Cells.Item(4,13).Value = 625    # M4 Precio promedio ponderado
$ws.Cells.Item(4,16).Value = 625    # P4 Precio $/Kg

# Row 5
$ws.Cells.Item(5,4).Value = 44608   # D5 Fecha
$ws.Cells.Item(5,10).Value = 120    # J5 Volumen (unchanged, kept for clarity)
$ws.Cells.Item(5,11).Value = 600    # K5 Precio minimo
$ws.Cells.Item(5,12).Value = 650    # L5 Precio maximo
$ws.Cells.Item(5,13).Value = 625    # M5 Precio promedio ponderado
$ws.Cells.Item(5,16).Value = 625    # P5 Precio $/Kg

# Row 6
$ws.Cells.Item(6,4).Value = 44532   # D6 Fecha
$ws.Cells.Item(6,10).Value = 60     # J6 Volumen

# Row 7
$ws.Cells.Item(7,4).Value = 44624   # D7 Fecha
$ws.Cells.Item(7,10).Value = 120    # J7 Volumen
$ws.Cells.Item(7,11).Value = 650    # K7 Precio minimo
$ws.Cells.Item(7,12).Value = 700    # L7 Precio maximo
$ws.Cells.Item(7,13).Value = 675    # M7 Precio promedio ponderado
$ws.Cells.Item(7,16).Value = 675    # P7 Precio $/Kg

# --- Append two new rows (8 and 9) with the older historical records ---

# Row 8
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8,3).Value = "Ñuble"
$ws.Cells.Item(8,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,4).Value = 44533
$ws.Cells.Item(8,5).Value = 16
$ws.Cells.Item(8,6).Value = 100112044
$ws.Cells.Item(8,7).Value = "Perejil"
$ws.Cells.Item(8,8).Value = "Sin especificar"
$ws.Cells.Item(8,9).Value = "Primera"
$ws.Cells.Item(8,10).Value = 100
$ws.Cells.Item(8,11).Value = 2000
$ws.Cells.Item(8,12).Value = 2200
$ws.Cells.Item(8,13).Value = 2100
$ws.Cells.Item(8,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8,15).Value = "Región del Maule"
$ws.Cells.Item(8,16).Value = 2100
$ws.Cells.Item(8,17).Value = 1
$ws.Cells.Item(8,18).Value = "Hortaliza"

# Row 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9,3).Value = "Ñuble"
$ws.Cells.Item(9,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,4).Value = 44754
$ws.Cells.Item(9,5).Value = 16
$ws.Cells.Item(9,6).Value = 100112044
$ws.Cells.Item(9,7).Value = "Perejil"
$ws.Cells.Item(9,8).Value = "Sin especificar"
$ws.Cells.Item(9,9).Value = "Primera"
$ws.Cells.Item(9,10).Value = 200
$ws.Cells.Item(9,11).Value = 700
$ws.Cells.Item(9,12).Value = 750
$ws.Cells.Item(9,13).Value = 725
$ws.Cells.Item(9,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9,15).Value = "Región del Maule"
$ws.Cells.Item(9,16).Value = 725
$ws.Cells.Item(9,17).Value = 1
$ws.Cells.Item(9,18).Value = "Hortaliza"

Write-Output "edit complete"
